$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet 1 (quality_comparison): add borders to C1/D1 header cells ---
# C1 -> top+bottom thin border (matches pre-existing border style "borderId=4")
$ws1.Range("C1").ClearFormats()
$ws1.Range("C1").Borders.Item(8).Weight = 2   # xlEdgeTop
$ws1.Range("C1").Borders.Item(9).Weight = 2   # xlEdgeBottom

# D1 -> top+right+bottom thin border (matches pre-existing border style "borderId=5")
$ws1.Range("D1").ClearFormats()
$ws1.Range("D1").Borders.Item(8).Weight = 2    # xlEdgeTop
$ws1.Range("D1").Borders.Item(10).Weight = 2   # xlEdgeRight
$ws1.Range("D1").Borders.Item(9).Weight = 2    # xlEdgeBottom

# Rename header label from "fedcore" to "approach"
$ws1.Range("C2").Value = "approach"

# --- Sheet 2 (computational_comparison): mirror the same border formatting ---
# Reuse the already-built styles from sheet1 (C1/D1) via copy/paste-format so
# the same cellXfs/border entries get reused instead of new ones being created.
$ws1.Range("C1").Copy()
$ws2.Range("C1").PasteSpecial(-4122)   # xlPasteFormats
$ws1.Range("D1").Copy()
$ws2.Range("D1").PasteSpecial(-4122)   # xlPasteFormats
$ws1.Range("C1").Copy()
$ws2.Range("F1").PasteSpecial(-4122)   # xlPasteFormats
$ws1.Range("D1").Copy()
$ws2.Range("G1").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# Rename header labels from "fedcore" to "approach"
$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

# Remove the stray empty inline-string cell G5
$ws2.Range("G5").ClearContents()
